# Slide 6 ("Алгоритм получения ассоциативных правил Apriori"):
# swap the two pictures' positions and re-stack the z-order so that
# "Рисунок 23" (id 24) moves behind the title/slide-number placeholders
# (sent to back) while "Рисунок 13" (id 14) moves above them (brought to
# the front). Sizes are unchanged; only Left/Top move.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

$pic13 = $s.Shapes.Item("Рисунок 13")
$pic23 = $s.Shapes.Item("Рисунок 23")

# msoZOrderCmd constants (not predefined in this host's PowerShell):
$msoSendToBack   = 1
$msoBringToFront = 0

# Reposition "Рисунок 23" (id 24) to its new spot (size unchanged).
# Target EMU: off x=6314075 y=3484072, ext unchanged (5548039 x 2798460).
$pic23.Left = 497.1712598425197
$pic23.Top  = 274.3363809527559

# Reposition "Рисунок 13" (id 14) to its new spot (size unchanged).
# Target EMU: off x=329886 y=1041003, ext unchanged (6522077 x 4128294).
$pic13.Left = 25.975276590551182
$pic13.Top  = 81.96874015748031

# Re-stack: send "Рисунок 23" all the way to the back, then bring
# "Рисунок 13" all the way to the front, so the final order (back -> front)
# is: Рисунок 23, Заголовок 1, Номер слайда 3, Рисунок 13.
$pic23.ZOrder($msoSendToBack)
$pic13.ZOrder($msoBringToFront)
